# Updates the cryptos list worksheet (price/volume figures, plus a few
# coin name/link/price swaps where source rows were reordered) to match
# the Wed Jun 26 07:57:10 UTC 2024 GitHub Actions data refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.687.51'
$ws.Range("E2").Value = '  +1.67%  '
$ws.Range("D3").Value = '3.399.81'
$ws.Range("E3").Value = '  +1.12%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = '''578.48'
$ws.Range("E5").Value = '  +1.27%  '
$ws.Range("D6").Value = '''137.92'
$ws.Range("E6").Value = '  +2.14%  '
$ws.Range("E7").Value = '  -0.08%  '
$ws.Range("D8").Value = '3.398.86'
$ws.Range("E8").Value = '  +1.11%  '
$ws.Range("E9").Value = '  -0.51%  '
$ws.Range("E10").Value = '  -0.33%  '
$ws.Range("E11").Value = '  +2.95%  '
$ws.Range("E12").Value = '  +0.18%  '
$ws.Range("D13").Value = '3.977.93'
$ws.Range("E13").Value = '  +1.13%  '
$ws.Range("D14").Value = '''0.124'
$ws.Range("E14").Value = '  +2.21%  '
$ws.Range("E15").Value = '  +3.13%  '
$ws.Range("D16").Value = '''25.99'
$ws.Range("E16").Value = '  +3.43%  '
$ws.Range("D17").Value = '3.398.72'
$ws.Range("E17").Value = '  +1.05%  '
$ws.Range("D18").Value = '61.753.53'
$ws.Range("E18").Value = '  +1.45%  '
$ws.Range("E19").Value = '  +2.56%  '
$ws.Range("D20").Value = '''5.88'
$ws.Range("D21").Value = '''9.46'
$ws.Range("E21").Value = '  +0.23%  '
$ws.Range("D22").Value = '''378.16'
$ws.Range("E22").Value = '  +1.77%  '
$ws.Range("E23").Value = '  -1.32%  '
$ws.Range("D24").Value = '3.530.12'
$ws.Range("E24").Value = '  +0.96%  '
$ws.Range("B25").Value = 'PEPE'
$ws.Range("C25").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D25").Value = '''0.0000128'
$ws.Range("E25").Value = '  +9.73%  '
$ws.Range("B26").Value = 'Dai'
$ws.Range("C26").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D26").Value = '''1.00'
$ws.Range("E26").Value = '  -0.08%  '
$ws.Range("D27").Value = '''71.32'
$ws.Range("E27").Value = '  +1.06%  '
$ws.Range("D28").Value = '''1.69'
$ws.Range("E28").Value = '  +0.30%  '
$ws.Range("E29").Value = '  -1.73%  '
$ws.Range("E30").Value = '  -0.20%  '
$ws.Range("B31").Value = 'InternetComputer(DFINITY)'
$ws.Range("C31").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D31").Value = '''8.26'
$ws.Range("E31").Value = '  +1.91%  '
$ws.Range("B32").Value = 'Kaspa'
$ws.Range("C32").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D32").Value = '''0.161'
$ws.Range("E32").Value = '  +4.07%  '
$ws.Range("E33").Value = '  +2.05%  '
$ws.Range("E34").Value = '  +0.04%  '
$ws.Range("E35").Value = '  +0.46%  '
$ws.Range("E36").Value = '  -3.29%  '
$ws.Range("E37").Value = '  +0.83%  '
$ws.Range("E38").Value = '  -0.54%  '
$ws.Range("D39").Value = '''165.46'
$ws.Range("E39").Value = '  +1.51%  '
$ws.Range("E40").Value = '  -0.17%  '
$ws.Range("E41").Value = '  +9.42%  '
$ws.Range("D42").Value = '''0.783'
$ws.Range("E42").Value = '  +3.28%  '
$ws.Range("B43").Value = 'ONDO'
$ws.Range("C43").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D43").Value = '''1.24'
$ws.Range("E43").Value = '  +3.21%  '
$ws.Range("B44").Value = 'FirstDigitalUSD'
$ws.Range("C44").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D44").Value = '''1.00'
$ws.Range("E44").Value = '  -0.06%  '
$ws.Range("D45").Value = '''25.23'
$ws.Range("E45").Value = '  +9.82%  '
$ws.Range("E46").Value = '  +0.84%  '
$ws.Range("D47").Value = '''41.53'
$ws.Range("E47").Value = '  +0.69%  '
$ws.Range("E48").Value = '  -1.51%  '
$ws.Range("D49").Value = '''22.84'
$ws.Range("E49").Value = '  -1.36%  '
$ws.Range("D50").Value = '2.342.13'
$ws.Range("E50").Value = '  +5.79%  '
$ws.Range("B51").Value = 'LidoDAOToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D51").Value = '''2.39'
$ws.Range("E51").Value = '  -1.49%  '
